# EDI/H_BL_NO generation algorithm update
#
# For every worksheet in the workbook, columns H ("gen_time") and I ("SHEET")
# swap places: H becomes "SHEET" (the sheet/tab name) and I becomes "gen_time".
# In addition, the gen_time value itself is refreshed from the old
# "2022-08-30 22:55" to the new "2022-09-01 23:53".

$wb = $excel.ActiveWorkbook
$newGenTime = "2022-09-01 23:53"

foreach ($ws in $wb.Worksheets) {
    $lastRow = $ws.UsedRange.Rows.Count

    # Swap the contents of column H (8) and column I (9) for every used row,
    # including the header row, so the header labels and every data value
    # trade places.
    for ($r = 1; $r -le $lastRow; $r++) {
        $hCell = $ws.Cells.Item($r, 8)
        $iCell = $ws.Cells.Item($r, 9)
        $hVal = $hCell.Value2
        $iVal = $iCell.Value2
        $hCell.Value = $iVal
        $iCell.Value = $hVal
    }

    # Column I now holds the gen_time value on every data row (row 1 is the
    # header). Refresh it to the new generation timestamp.
    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 9).Value = $newGenTime
    }

    # Keep the column widths in sync with the swapped content.
    $hCol = $ws.Columns.Item(8)
    $iCol = $ws.Columns.Item(9)
    $hWidth = $hCol.ColumnWidth
    $iWidth = $iCol.ColumnWidth
    $hCol.ColumnWidth = $iWidth
    $iCol.ColumnWidth = $hWidth
}
